$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D7").Value = "2016-03-04 10:37:23"
$wsZh.Range("G7").Value = "2016-03-04 10:38:39"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D7").Value = "2016-03-04 10:37:38"
$wsDe.Range("G7").Value = "2016-03-04 10:39:16"
